$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then numeric columns E..T
$rows = @(
    @("ECs","Dlk1","Notch4","ECs",   1, 0.3333333333333333, 0.09821600000000001, 0.294648, 0.02540301315491012, 0.02540301315491011, 3, 1, 29.546424, 88.63927200000001, 0.9033225104610835, 0.9033225104610834, 2.901931579584001, 26.117384216256, 0.02294711361636934, 0.02294711361636933),
    @("ECs","Dlk1","Notch4","FAPs",  1, 0.3333333333333333, 0.09821600000000001, 0.294648, 0.02540301315491012, 0.02540301315491011, 3, 1, 1.672785333333334, 5.018356000000001, 0.05114204841740398, 0.05114204841740398, 0.1642942842986667, 1.478648558688, 0.001299162128716364, 0.001299162128716363),
    @("ECs","Dlk1","Notch4","sCs",   1, 0.3333333333333333, 0.09821600000000001, 0.294648, 0.02540301315491012, 0.02540301315491011, 3, 1, 1.489401, 4.468203, 0.04553544112151264, 0.04553544112151264, 0.146283008616, 1.316547077544, 0.001156737409824421, 0.00115673740982442),
    @("FAPs","Dlk1","Notch4","ECs",  3, 1, 2.101448, 6.304344, 0.5435276450716743, 0.5435276450716742, 3, 1, 29.546424, 88.63927200000001, 0.9033225104610835, 0.9033225104610834, 62.090273621952, 558.812462597568, 0.4909807568511456, 0.4909807568511454),
    @("FAPs","Dlk1","Notch4","FAPs", 3, 1, 2.101448, 6.304344, 0.5435276450716743, 0.5435276450716742, 3, 1, 1.672785333333334, 5.018356000000001, 0.05114204841740398, 0.05114204841740398, 3.515271393162667, 31.63744253846401, 0.02779711714045313, 0.02779711714045313),
    @("FAPs","Dlk1","Notch4","sCs",  3, 1, 2.101448, 6.304344, 0.5435276450716743, 0.5435276450716742, 3, 1, 1.489401, 4.468203, 0.04553544112151264, 0.04553544112151264, 3.129898752648, 28.169088773832, 0.02474977108007564, 0.02474977108007564),
    @("sCs","Dlk1","Notch4","ECs",   3, 1, 1.666649, 4.999947, 0.4310693417734156, 0.4310693417734155, 3, 1, 29.546424, 88.63927200000001, 0.9033225104610835, 0.9033225104610834, 49.243518013176, 443.191662118584, 0.3893946399935686, 0.3893946399935685),
    @("sCs","Dlk1","Notch4","FAPs",  3, 1, 1.666649, 4.999947, 0.4310693417734156, 0.4310693417734155, 3, 1, 1.672785333333334, 5.018356000000001, 0.05114204841740398, 0.05114204841740398, 2.787946003014667, 25.091514027132, 0.02204576914823449, 0.02204576914823448),
    @("sCs","Dlk1","Notch4","sCs",   3, 1, 1.666649, 4.999947, 0.4310693417734156, 0.4310693417734155, 3, 1, 1.489401, 4.468203, 0.04553544112151264, 0.04553544112151264, 2.482308687249, 22.340778185241, 0.01962893263161258, 0.01962893263161257)
)

$rowIndex = 2
foreach ($row in $rows) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
